# Fix typo on the "Project outputs" slide: "the 5 csv files" -> "the 4 csv files"
# (slide 6, placeholder body text box that starts with "Send the 5 csv files + ...")

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -like "*the 5 csv files*") {
                $targetShape = $shp
                $targetSlide = $s
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the exact substring "the 5 " (including the trailing space) and
# replace it with "the 4 " so the run is split the same way PowerPoint
# splits a run when only part of its text is edited in place.
$hit = $tr.Find("the 5 ", 0, 0)
$hit.Text = "the 4 "
